$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$range = $ws.Range("A1:U82")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.TableStyle = "TableStyleLight1"
Write-Output "ok"
